# Add a new "logout_link" locator row beneath the existing table and
# move the selection to the next empty row, mirroring how a user would
# continue filling in the locator sheet after parsing new command-line
# arguments for the logout element.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "logout_link"
$ws.Range("B6").Value = "link text"
$ws.Range("C6").Value = "Log out"

$ws.Range("A7").Select()
